$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New result columns for block sizes 256 and 512, added alongside the
# existing 32 / 64 / 128 columns (A:D stay untouched).
$ws.Range("E1").Value = "'256"
$ws.Range("F1").Value = "'512"

# Match the header formatting (bold, bordered, centered) used by A1:D1.
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Per-method timing values for the 256 and 512 columns, row-aligned with
# the existing data in rows 2-65.
$data = @(
    @(953, 1982),
    @(953, 1982),
    @(1386, 2817),
    @(1386, 2817),
    @(1674, 3398),
    @(1674, 3398),
    @(1209, 2700),
    @(1209, 2700),
    @(661, 1186),
    @(661, 1186),
    @(2359, 4641),
    @(2359, 4641),
    @(1838, 3769),
    @(1838, 3769),
    @(1438, 2824),
    @(1438, 2824),
    @(2294, 4775),
    @(2294, 4775),
    @(1308, 2589),
    @(1308, 2589),
    @(2311, 4819),
    @(2311, 4819),
    @(1599, 3301),
    @(1599, 3301),
    @(1341, 2520),
    @(1341, 2520),
    @(2564, 6900),
    @(2564, 6900),
    @(2515, 4584),
    @(2515, 4584),
    @(1183, 2656),
    @(1183, 2656),
    @(16, 17),
    @(16, 17),
    @(24, 24),
    @(24, 24),
    @(15, 15),
    @(15, 15),
    @(18, 19),
    @(18, 19),
    @(30, 27),
    @(30, 27),
    @(29, 27),
    @(29, 27),
    @(31, 29),
    @(31, 29),
    @(28, 28),
    @(28, 28),
    @(12, 12),
    @(12, 12),
    @(24, 24),
    @(24, 24),
    @(15, 15),
    @(15, 15),
    @(19, 19),
    @(19, 19),
    @(29, 29),
    @(29, 29),
    @(31, 29),
    @(31, 29),
    @(33, 29),
    @(33, 29),
    @(33, 29),
    @(33, 29)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $data[$i][0]
    $ws.Cells.Item($row, 6).Value = $data[$i][1]
}
